$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.934.77'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '1.773.12'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4501'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3579'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07479'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.20'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.099'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.064'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.223'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").Value = '1.773.79'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001063'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06445'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.834'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").Value = '27.946.04'
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.113'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").Value = '1.972.50'
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.211'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.106'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09168'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.582'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.635'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02302'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06122'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2100'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6362'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.985'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.187'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.395'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.959'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.740'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5897'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.964'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.143'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06927'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.34%  '
